$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 21.65036466666666
$ws.Range("H2").Value = 64.951094
$ws.Range("I2").Value = 0.9284333993050746
$ws.Range("J2").Value = 0.9284333993050747
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.06661333333333334
$ws.Range("N2").Value = 0.19984
$ws.Range("O2").Value = 0.01028428344934373
$ws.Range("P2").Value = 0.01028428344934373
$ws.Range("Q2").Value = 1.442202958328889
$ws.Range("R2").Value = 12.97982662496
$ws.Range("S2").Value = 0.009548272242291116
$ws.Range("T2").Value = 0.009548272242291119
$ws.Range("G3").Value = 21.65036466666666
$ws.Range("H3").Value = 64.951094
$ws.Range("I3").Value = 0.9284333993050746
$ws.Range("J3").Value = 0.9284333993050747
$ws.Range("O3").Value = 0.7439394054794892
$ws.Range("P3").Value = 0.7439394054794893
$ws.Range("Q3").Value = 104.3253637149042
$ws.Range("R3").Value = 938.9282734341379
$ws.Range("S3").Value = 0.6906981911063185
$ws.Range("T3").Value = 0.6906981911063186
$ws.Range("G4").Value = 21.65036466666666
$ws.Range("H4").Value = 64.951094
$ws.Range("I4").Value = 0.9284333993050746
$ws.Range("J4").Value = 0.9284333993050747
$ws.Range("O4").Value = 0.245776311071167
$ws.Range("P4").Value = 0.245776311071167
$ws.Range("Q4").Value = 34.46611761139444
$ws.Range("R4").Value = 310.19505850255
$ws.Range("S4").Value = 0.228186935956465
$ws.Range("T4").Value = 0.2281869359564651
$ws.Range("G5").Value = 0.96805
$ws.Range("I5").Value = 0.04151292442575075
$ws.Range("J5").Value = 0.04151292442575075
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.06661333333333334
$ws.Range("N5").Value = 0.19984
$ws.Range("O5").Value = 0.01028428344934373
$ws.Range("P5").Value = 0.01028428344934373
$ws.Range("Q5").Value = 0.06448503733333334
$ws.Range("R5").Value = 0.5803653360000001
$ws.Range("S5").Value = 0.0004269306816056054
$ws.Range("T5").Value = 0.0004269306816056055
$ws.Range("G6").Value = 0.96805
$ws.Range("I6").Value = 0.04151292442575075
$ws.Range("J6").Value = 0.04151292442575075
$ws.Range("O6").Value = 0.7439394054794892
$ws.Range("P6").Value = 0.7439394054794893
$ws.Range("Q6").Value = 4.664686710783333
$ws.Range("R6").Value = 41.98218039704999
$ws.Range("S6").Value = 0.03088310031700798
$ws.Range("T6").Value = 0.03088310031700798
$ws.Range("G7").Value = 0.96805
$ws.Range("I7").Value = 0.04151292442575075
$ws.Range("J7").Value = 0.04151292442575075
$ws.Range("O7").Value = 0.245776311071167
$ws.Range("P7").Value = 0.245776311071167
$ws.Range("S7").Value = 0.01020289342713716
$ws.Range("T7").Value = 0.01020289342713716
$ws.Range("I8").Value = 0.0300536762691746
$ws.Range("J8").Value = 0.0300536762691746
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.06661333333333334
$ws.Range("N8").Value = 0.19984
$ws.Range("O8").Value = 0.01028428344934373
$ws.Range("P8").Value = 0.01028428344934373
$ws.Range("Q8").Value = 0.04668455578666668
$ws.Range("R8").Value = 0.42016100208
$ws.Range("S8").Value = 0.0003090805254470067
$ws.Range("T8").Value = 0.0003090805254470067
$ws.Range("I9").Value = 0.0300536762691746
$ws.Range("J9").Value = 0.0300536762691746
$ws.Range("O9").Value = 0.7439394054794892
$ws.Range("P9").Value = 0.7439394054794893
$ws.Range("S9").Value = 0.02235811405616278
$ws.Range("T9").Value = 0.02235811405616279
$ws.Range("I10").Value = 0.0300536762691746
$ws.Range("J10").Value = 0.0300536762691746
$ws.Range("O10").Value = 0.245776311071167
$ws.Range("P10").Value = 0.245776311071167
$ws.Range("S10").Value = 0.007386481687564807
$ws.Range("T10").Value = 0.007386481687564806
